$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert the worksheet's data range into an Excel Table (ListObject),
# matching the author's "Converted worksheets into excel table" commit.
# The table is created over the full A:D columns (as Excel does when the
# entire columns are selected before Insert > Table), which yields the
# same ref="A1:D1048576" seen in the target workbook.
$rng = $ws.Range("A1:D1048576")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = "TableStyleMedium2"

# Widen column B (HighestInfectionCount) to fit the new table header/filter
# button, as Excel does automatically when a table is created.
$ws.Columns.Item(2).ColumnWidth = 22.29

# Restore the view state recorded after the edit: the window had been
# scrolled down and a different cell was selected.
$excel.ActiveWindow.ScrollRow = 49
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("H9").Select()
